# Atualização automática de preços de eletricidade
# Updates row 2 of the SpotPTTable with the latest daily hourly spot prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45993
$ws.Range("B2").Value = 79.03
$ws.Range("C2").Value = 76.36
$ws.Range("D2").Value = 72.23
$ws.Range("E2").Value = 68.77
$ws.Range("F2").Value = 67.95999999999999
$ws.Range("G2").Value = 70.09999999999999
$ws.Range("H2").Value = 81
$ws.Range("I2").Value = 90.29000000000001
$ws.Range("J2").Value = 95.48999999999999
$ws.Range("K2").Value = 91.42
$ws.Range("L2").Value = 84.89
$ws.Range("M2").Value = 82.51000000000001
$ws.Range("N2").Value = 80.3
$ws.Range("O2").Value = 73.90000000000001
$ws.Range("P2").Value = 72.42
$ws.Range("Q2").Value = 84.45999999999999
$ws.Range("R2").Value = 87.68000000000001
$ws.Range("S2").Value = 111.43
$ws.Range("T2").Value = 125.89
$ws.Range("U2").Value = 126.32
$ws.Range("V2").Value = 120.3
$ws.Range("W2").Value = 103.16
$ws.Range("X2").Value = 94.03
$ws.Range("Y2").Value = 87.68000000000001
$ws.Range("Z2").Value = 88.65000000000001

$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 112.83
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 126.1
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 111.73
$ws.Range("AG2").Value = "0h-23h"
